$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 85 appended to the "jogos" table (Nome do Jogo / Status / Plataforma / Objetivo).
# Column A holds "3" which must be stored as text (like the "1"/"2" values already present
# in A83/A84), so the cell is formatted as Text before the value is written - otherwise a
# numeric-looking string like "3" would be auto-converted to a real number by Excel.
$ws.Range("A85").NumberFormat = "@"
$ws.Range("A85").Value = "3"
$ws.Range("B85").Value = "Incompleto"
$ws.Range("C85").Value = "PS3"
$ws.Range("D85").Value = "Concluído"
